# Weekly re-sync of daily price rows for "Fruta, Feria Lagunitas de Puerto Montt - Granada"
# (commit: "Fruta / hortaliza, semanal") -- updates Fecha/Calidad/Volumen/Precio* columns
# on the existing rows to reflect the weekly-aggregated source values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @{
    2 = @{ "D" = 44351; "M" = 100; "N" = 15000; "O" = 15000; "P" = 15000; "S" = 1000 }
    3 = @{ "D" = 44351; "M" = 200; "N" = 13000; "O" = 13500; "P" = 13250; "S" = 883 }
    4 = @{ "D" = 44348; "M" = 100; "N" = 15000; "O" = 15000; "P" = 15000; "S" = 1000 }
    5 = @{ "D" = 44348; "M" = 200; "N" = 13000; "O" = 13500; "P" = 13250; "S" = 883 }
    6 = @{ "D" = 44323; "M" = 100; "N" = 17000; "O" = 17000; "P" = 17000; "S" = 1133 }
    7 = @{ "D" = 44323; "M" = 100; "O" = 14000; "P" = 14000; "S" = 933 }
    8 = @{ "D" = 44344; "M" = 100; "N" = 16000; "O" = 16000; "P" = 16000; "S" = 1067 }
    9 = @{ "D" = 44344; "L" = "Segunda"; "M" = 120; "N" = 13000; "O" = 13500; "P" = 13250; "S" = 883 }
    10 = @{ "D" = 44285; "L" = "Segunda"; "M" = 160; "N" = 15000; "O" = 16000; "P" = 15500; "S" = 1033 }
    11 = @{ "L" = "Primera"; "M" = 60; "N" = 17000; "O" = 17000; "P" = 17000; "S" = 1133 }
    12 = @{ "D" = 44336; "L" = "Segunda"; "M" = 120; "N" = 14000; "O" = 14500; "P" = 14250; "S" = 950 }
    13 = @{ "D" = 44302; "L" = "Primera"; "M" = 100; "N" = 17500; "O" = 17500; "P" = 17500; "S" = 1167 }
    14 = @{ "D" = 44302; "L" = "Segunda"; "M" = 200; "N" = 14000; "P" = 14500; "S" = 967 }
    15 = @{ "D" = 44309; "L" = "Primera"; "M" = 100; "N" = 17500; "O" = 17500; "P" = 17500; "S" = 1167 }
    16 = @{ "L" = "Segunda"; "M" = 200; "N" = 14000; "O" = 14500; "P" = 14250; "S" = 950 }
    17 = @{ "D" = 44292; "M" = 160; "S" = 967 }
    18 = @{ "D" = 44299; "L" = "Primera"; "M" = 60; "N" = 17500; "O" = 17500; "P" = 17500; "S" = 1167 }
    19 = @{ "L" = "Segunda"; "M" = 120; "N" = 14000; "O" = 15000; "P" = 14500; "S" = 967 }
    20 = @{ "D" = 44306; "L" = "Primera"; "M" = 100; "N" = 17500; "O" = 17500; "P" = 17500; "S" = 1167 }
    21 = @{ "D" = 44306; "L" = "Segunda"; "M" = 200; "N" = 14000; "O" = 14500; "P" = 14250; "S" = 950 }
    22 = @{ "D" = 44295; "L" = "Primera"; "M" = 160; "S" = 967 }
    25 = @{ "D" = 44313; "L" = "Especial"; "N" = 17500; "O" = 17500; "P" = 17500; "Q" = "$/caja 14 kilos empedrada"; "S" = 1250; "T" = 14 }
    26 = @{ "D" = 44313; "L" = "Primera"; "M" = 100; "N" = 16000; "O" = 16000; "P" = 16000; "Q" = "$/caja 14 kilos empedrada"; "S" = 1143; "T" = 14 }
    27 = @{ "D" = 44313; "L" = "Segunda"; "M" = 80; "N" = 14000; "O" = 14000; "P" = 14000; "Q" = "$/caja 14 kilos empedrada"; "S" = 1000; "T" = 14 }
    28 = @{ "D" = 44334; "L" = "Primera"; "M" = 200; "O" = 17000; "P" = 15500; "S" = 1033 }
    29 = @{ "D" = 44334; "L" = "Segunda"; "N" = 14500; "O" = 14500; "P" = 14500; "S" = 967 }
    30 = @{ "D" = 44330; "L" = "Primera"; "M" = 100; "N" = 17000; "O" = 17000; "P" = 17000; "S" = 1133 }
    31 = @{ "D" = 44330; "M" = 200; "N" = 14000; "O" = 14500; "P" = 14250; "S" = 950 }
    33 = @{ "D" = 44305; "M" = 60 }
    34 = @{ "D" = 44305; "M" = 120 }
    35 = @{ "D" = 44301; "M" = 60; "N" = 17500; "O" = 17500; "P" = 17500; "S" = 1167 }
    36 = @{ "D" = 44301; "M" = 80; "N" = 14000; "O" = 15000; "P" = 14500; "S" = 967 }
    37 = @{ "D" = 44293; "L" = "Primera"; "M" = 60; "N" = 14000; "O" = 15000; "P" = 14500; "Q" = "$/caja 15 kilos empedrada"; "S" = 967; "T" = 15 }
    38 = @{ "D" = 44327; "N" = 17000; "O" = 17000; "P" = 17000; "Q" = "$/caja 15 kilos empedrada"; "S" = 1133; "T" = 15 }
    39 = @{ "D" = 44327; "M" = 200; "O" = 14500; "P" = 14250; "Q" = "$/caja 15 kilos empedrada"; "S" = 950; "T" = 15 }
}

foreach ($row in $rowUpdates.Keys) {
    foreach ($col in $rowUpdates[$row].Keys) {
        $ws.Range("$col$row").Value = $rowUpdates[$row][$col]
    }
}